$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93; this shifts the existing rows 93:165 down to 94:166
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new record's data
$ws.Cells.Item(93, 1).Value = 10
$ws.Cells.Item(93, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(93, 3).Value = "La Araucanía"
$ws.Cells.Item(93, 4).Value = 44904
$ws.Cells.Item(93, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(93, 5).Value = 9
$ws.Cells.Item(93, 6).Value = 100112031
$ws.Cells.Item(93, 7).Value = "Poroto verde"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 135
$ws.Cells.Item(93, 11).Value = 35000
$ws.Cells.Item(93, 12).Value = 37000
$ws.Cells.Item(93, 13).Value = 36037
$ws.Cells.Item(93, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(93, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(93, 16).Value = 1441
$ws.Cells.Item(93, 17).Value = 25
$ws.Cells.Item(93, 18).Value = "Hortaliza"
